$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.837.51"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = "'3.393.77"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'569.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').Value = "'162.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'3.391.40"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').Value = "'0.545"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.10%  '
$ws.Range('D10').Value = "'7.29"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('E11').Value = '  -2.58%  '
$ws.Range('D12').Value = "'0.419"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.64%  '
$ws.Range('D13').Value = "'3.982.60"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = "'26.87"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('D17').Value = "'63.894.49"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = "'3.399.23"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').Value = "'6.09"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('D20').Value = "'13.51"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').Value = "'376.09"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = "'7.75"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = "'70.10"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.05%  '
$ws.Range('E25').Value = '  -4.69%  '
$ws.Range('D26').Value = "'0.0000113"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.23%  '
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('D28').Value = "'0.179"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = "'6.06"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = "'1.38"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.22%  '
$ws.Range('D32').Value = "'1.99"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D34').Value = "'22.74"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').Value = "'6.98"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -5.90%  '
$ws.Range('D37').Value = "'159.55"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('E38').Value = '  +9.29%  '
$ws.Range('D39').Value = "'1.80"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('E40').Value = '  -3.86%  '
$ws.Range('D41').Value = "'25.70"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').Value = "'42.71"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.72%  '
$ws.Range('D43').Value = "'2.732.90"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.78%  '
$ws.Range('D44').Value = "'26.06"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = "'6.40"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.73%  '
$ws.Range('D46').Value = "'4.35"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('D47').Value = "'0.0304"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('D48').Value = "'2.39"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').Value = "'327.77"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('E50').Value = '  -5.05%  '
$ws.Range('E51').Value = '  -1.99%  '
